$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 11
$ws.Range("C4").Value = 12
$ws.Range("C5").Value = 7
$ws.Range("C6").Value = 9
$ws.Range("C7").Value = 5
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 2
$ws.Range("C10").Value = 7
$ws.Range("C11").Value = 7
$ws.Range("C12").Value = 5
$ws.Range("C14").Value = 6
$ws.Range("C15").Value = 9
$ws.Range("C16").Value = 6
$ws.Range("C18").Value = 8
